$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.350.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.476.60"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.28"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.10"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.53%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.477.10"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.84%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.05%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.19"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +9.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000180"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.918.57"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.245.84"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.460.96"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.11"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "331.82"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.24"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +10.08%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "678.22"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +8.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "66.33"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.47%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +9.38%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.44%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.38%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.20"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.88"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.138"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.997"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.82%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "154.12"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.375"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.91"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.79"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.52%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.10%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0301"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.73%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.15"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +27.67%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.11"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.37%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.64"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.99"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.65%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.610"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("B51").Value = "Hedera"
$ws.Range("C51").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0519"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.90%  "
